$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The "From" value for rule R30 (row 10) changed from 18 to 1.
$ws.Range("C10").Value = 1

